$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Timp3"
$ws.Cells.Item(2,3).Value = "Kdr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 113.9271086666667
$ws.Cells.Item(2,8).Value = 341.781326
$ws.Cells.Item(2,9).Value = 0.4186548232357613
$ws.Cells.Item(2,10).Value = 0.5001257596740779
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 153.3847273333333
$ws.Cells.Item(2,14).Value = 460.154182
$ws.Cells.Item(2,15).Value = 0.9432535557163702
$ws.Cells.Item(2,16).Value = 0.9461442014370421
$ws.Cells.Item(2,17).Value = 17474.6784987117
$ws.Cells.Item(2,18).Value = 157272.1064884053
$ws.Cells.Item(2,19).Value = 0.3948976506349403
$ws.Cells.Item(2,20).Value = 0.4731910875049245
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Timp3"
$ws.Cells.Item(3,3).Value = "Kdr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 113.9271086666667
$ws.Cells.Item(3,8).Value = 341.781326
$ws.Cells.Item(3,9).Value = 0.4186548232357613
$ws.Cells.Item(3,10).Value = 0.5001257596740779
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.06306133333333333
$ws.Cells.Item(3,14).Value = 0.189184
$ws.Cells.Item(3,15).Value = 0.0003878014971178633
$ws.Cells.Item(3,16).Value = 0.0003889899333885992
$ws.Cells.Item(3,17).Value = 7.184395375331555
$ws.Cells.Item(3,18).Value = 64.659558377984
$ws.Cells.Item(3,19).Value = 0.0001623549672264427
$ws.Cells.Item(3,20).Value = 0.0001945438859415421
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Timp3"
$ws.Cells.Item(4,3).Value = "Kdr"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 113.9271086666667
$ws.Cells.Item(4,8).Value = 341.781326
$ws.Cells.Item(4,9).Value = 0.4186548232357613
$ws.Cells.Item(4,10).Value = 0.5001257596740779
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.722979333333332
$ws.Cells.Item(4,14).Value = 17.168938
$ws.Cells.Item(4,15).Value = 0.0351939902968738
$ws.Cells.Item(4,16).Value = 0.03530184396657746
$ws.Cells.Item(4,17).Value = 652.002488405754
$ws.Cells.Item(4,18).Value = 5868.022395651787
$ws.Cells.Item(4,19).Value = 0.0147341337866988
$ws.Cells.Item(4,20).Value = 0.01765536153168032
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Timp3"
$ws.Cells.Item(5,3).Value = "Kdr"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 113.9271086666667
$ws.Cells.Item(5,8).Value = 341.781326
$ws.Cells.Item(5,9).Value = 0.4186548232357613
$ws.Cells.Item(5,10).Value = 0.5001257596740779
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.951202
$ws.Cells.Item(5,14).Value = 5.853605999999999
$ws.Cells.Item(5,15).Value = 0.01199909701844822
$ws.Cells.Item(5,16).Value = 0.0120358688262385
$ws.Cells.Item(5,17).Value = 222.2948022846173
$ws.Cells.Item(5,18).Value = 2000.653220561556
$ws.Cells.Item(5,19).Value = 0.005023479841247188
$ws.Cells.Item(5,20).Value = 0.006019448040060083
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Timp3"
$ws.Cells.Item(6,3).Value = "Kdr"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 113.9271086666667
$ws.Cells.Item(6,8).Value = 341.781326
$ws.Cells.Item(6,9).Value = 0.4186548232357613
$ws.Cells.Item(6,10).Value = 0.5001257596740779
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.490433
$ws.Cells.Item(6,14).Value = 2.980866
$ws.Cells.Item(6,15).Value = 0.009165555471189982
$ws.Cells.Item(6,16).Value = 0.006129095836753322
$ws.Cells.Item(6,17).Value = 169.800722351386
$ws.Cells.Item(6,18).Value = 1018.804334108316
$ws.Cells.Item(6,19).Value = 0.003837204005648607
$ws.Cells.Item(6,20).Value = 0.003065318711471483
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Timp3"
$ws.Cells.Item(7,3).Value = "Kdr"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 25.210481
$ws.Cells.Item(7,8).Value = 75.631443
$ws.Cells.Item(7,9).Value = 0.09264247631899748
$ws.Cells.Item(7,10).Value = 0.110670858845055
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 153.3847273333333
$ws.Cells.Item(7,14).Value = 460.154182
$ws.Cells.Item(7,15).Value = 0.9432535557163702
$ws.Cells.Item(7,16).Value = 0.9461442014370421
$ws.Cells.Item(7,17).Value = 3866.902754127181
$ws.Cells.Item(7,18).Value = 34802.12478714462
$ws.Cells.Item(7,19).Value = 0.08738534519826399
$ws.Cells.Item(7,20).Value = 0.1047105913643062
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Timp3"
$ws.Cells.Item(8,3).Value = "Kdr"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 25.210481
$ws.Cells.Item(8,8).Value = 75.631443
$ws.Cells.Item(8,9).Value = 0.09264247631899748
$ws.Cells.Item(8,10).Value = 0.110670858845055
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.06306133333333333
$ws.Cells.Item(8,14).Value = 0.189184
$ws.Cells.Item(8,15).Value = 0.0003878014971178633
$ws.Cells.Item(8,16).Value = 0.0003889899333885992
$ws.Cells.Item(8,17).Value = 1.589806545834667
$ws.Cells.Item(8,18).Value = 14.308258912512
$ws.Cells.Item(8,19).Value = 0.00003592689101321342
$ws.Cells.Item(8,20).Value = 0.00004304985001019701
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Timp3"
$ws.Cells.Item(9,3).Value = "Kdr"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 25.210481
$ws.Cells.Item(9,8).Value = 75.631443
$ws.Cells.Item(9,9).Value = 0.09264247631899748
$ws.Cells.Item(9,10).Value = 0.110670858845055
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 5.722979333333332
$ws.Cells.Item(9,14).Value = 17.168938
$ws.Cells.Item(9,15).Value = 0.0351939902968738
$ws.Cells.Item(9,16).Value = 0.03530184396657746
$ws.Cells.Item(9,17).Value = 144.2790617463926
$ws.Cells.Item(9,18).Value = 1298.511555717534
$ws.Cells.Item(9,19).Value = 0.003260458412649158
$ws.Cells.Item(9,20).Value = 0.00390688539059525
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Timp3"
$ws.Cells.Item(10,3).Value = "Kdr"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 25.210481
$ws.Cells.Item(10,8).Value = 75.631443
$ws.Cells.Item(10,9).Value = 0.09264247631899748
$ws.Cells.Item(10,10).Value = 0.110670858845055
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.951202
$ws.Cells.Item(10,14).Value = 5.853605999999999
$ws.Cells.Item(10,15).Value = 0.01199909701844822
$ws.Cells.Item(10,16).Value = 0.0120358688262385
$ws.Cells.Item(10,17).Value = 49.190740948162
$ws.Cells.Item(10,18).Value = 442.7166685334579
$ws.Cells.Item(10,19).Value = 0.001111626061380942
$ws.Cells.Item(10,20).Value = 0.001332019939946239
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Timp3"
$ws.Cells.Item(11,3).Value = "Kdr"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 25.210481
$ws.Cells.Item(11,8).Value = 75.631443
$ws.Cells.Item(11,9).Value = 0.09264247631899748
$ws.Cells.Item(11,10).Value = 0.110670858845055
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 1.490433
$ws.Cells.Item(11,14).Value = 2.980866
$ws.Cells.Item(11,15).Value = 0.009165555471189982
$ws.Cells.Item(11,16).Value = 0.006129095836753322
$ws.Cells.Item(11,17).Value = 37.57453282827301
$ws.Cells.Item(11,18).Value = 225.447196969638
$ws.Cells.Item(11,19).Value = 0.0008491197556901757
$ws.Cells.Item(11,20).Value = 0.0006783123001971412
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Timp3"
$ws.Cells.Item(12,3).Value = "Kdr"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 132.9889985
$ws.Cells.Item(12,8).Value = 265.977997
$ws.Cells.Item(12,9).Value = 0.4887027004452411
$ws.Cells.Item(12,10).Value = 0.3892033814808671
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 153.3847273333333
$ws.Cells.Item(12,14).Value = 460.154182
$ws.Cells.Item(12,15).Value = 0.9432535557163702
$ws.Cells.Item(12,16).Value = 0.9461442014370421
$ws.Cells.Item(12,17).Value = 20398.48127325558
$ws.Cells.Item(12,18).Value = 122390.8876395335
$ws.Cells.Item(12,19).Value = 0.4609705598831658
$ws.Cells.Item(12,20).Value = 0.3682425225678115
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Timp3"
$ws.Cells.Item(13,3).Value = "Kdr"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 132.9889985
$ws.Cells.Item(13,8).Value = 265.977997
$ws.Cells.Item(13,9).Value = 0.4887027004452411
$ws.Cells.Item(13,10).Value = 0.3892033814808671
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.06306133333333333
$ws.Cells.Item(13,14).Value = 0.189184
$ws.Cells.Item(13,15).Value = 0.0003878014971178633
$ws.Cells.Item(13,16).Value = 0.0003889899333885992
$ws.Cells.Item(13,17).Value = 8.386463564074667
$ws.Cells.Item(13,18).Value = 50.318781384448
$ws.Cells.Item(13,19).Value = 0.0001895196388782072
$ws.Cells.Item(13,20).Value = 0.0001513961974368601
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Timp3"
$ws.Cells.Item(14,3).Value = "Kdr"
$ws.Cells.Item(14,4).Value = "M1"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 132.9889985
$ws.Cells.Item(14,8).Value = 265.977997
$ws.Cells.Item(14,9).Value = 0.4887027004452411
$ws.Cells.Item(14,10).Value = 0.3892033814808671
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 5.722979333333332
$ws.Cells.Item(14,14).Value = 17.168938
$ws.Cells.Item(14,15).Value = 0.0351939902968738
$ws.Cells.Item(14,16).Value = 0.03530184396657746
$ws.Cells.Item(14,17).Value = 761.0932899761975
$ws.Cells.Item(14,18).Value = 4566.559739857185
$ws.Cells.Item(14,19).Value = 0.01719939809752584
$ws.Cells.Item(14,20).Value = 0.01373959704430189
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Timp3"
$ws.Cells.Item(15,3).Value = "Kdr"
$ws.Cells.Item(15,4).Value = "M2"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 132.9889985
$ws.Cells.Item(15,8).Value = 265.977997
$ws.Cells.Item(15,9).Value = 0.4887027004452411
$ws.Cells.Item(15,10).Value = 0.3892033814808671
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 1.951202
$ws.Cells.Item(15,14).Value = 5.853605999999999
$ws.Cells.Item(15,15).Value = 0.01199909701844822
$ws.Cells.Item(15,16).Value = 0.0120358688262385
$ws.Cells.Item(15,17).Value = 259.488399851197
$ws.Cells.Item(15,18).Value = 1556.930399107182
$ws.Cells.Item(15,19).Value = 0.005863991115820084
$ws.Cells.Item(15,20).Value = 0.00468440084623218
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Timp3"
$ws.Cells.Item(16,3).Value = "Kdr"
$ws.Cells.Item(16,4).Value = "sCs"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 132.9889985
$ws.Cells.Item(16,8).Value = 265.977997
$ws.Cells.Item(16,9).Value = 0.4887027004452411
$ws.Cells.Item(16,10).Value = 0.3892033814808671
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 1.490433
$ws.Cells.Item(16,14).Value = 2.980866
$ws.Cells.Item(16,15).Value = 0.009165555471189982
$ws.Cells.Item(16,16).Value = 0.006129095836753322
$ws.Cells.Item(16,17).Value = 198.2111920013505
$ws.Cells.Item(16,18).Value = 792.8447680054021
$ws.Cells.Item(16,19).Value = 0.004479231709851198
$ws.Cells.Item(16,20).Value = 0.002385464825084698